# Insert a new data row at row 38, pushing existing rows 38-58 down to 39-59.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(38).Insert()

# Fill in the new row 38 with the new record's values.
$ws.Cells.Item(38, 1).Value  = 1
$ws.Cells.Item(38, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(38, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(38, 4).Value  = 44629
$ws.Cells.Item(38, 5).Value  = 15
$ws.Cells.Item(38, 6).Value  = 100112009
$ws.Cells.Item(38, 7).Value  = "Acelga"
$ws.Cells.Item(38, 8).Value  = "Sin especificar"
$ws.Cells.Item(38, 9).Value  = "Primera"
$ws.Cells.Item(38, 10).Value = 270
$ws.Cells.Item(38, 11).Value = 1000
$ws.Cells.Item(38, 12).Value = 1200
$ws.Cells.Item(38, 13).Value = 1100
$ws.Cells.Item(38, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(38, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(38, 16).Value = 367
$ws.Cells.Item(38, 17).Value = 3
$ws.Cells.Item(38, 18).Value = "Hortaliza"
